# This script re-orders the "D" through "T" data for rows 2-8 and 10-13
# (row 9 is untouched) according to the target diff. The underlying data
# values themselves are unchanged; only which row they land on changes
# (a permutation of existing rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that participate in the move (D, K, L, M, N, O, P, Q, R, S, T).
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Source row -> destination row mapping (row 9 stays put).
$rowMap = @{
    2  = 11
    3  = 7
    4  = 8
    5  = 4
    6  = 10
    7  = 5
    8  = 6
    10 = 13
    11 = 2
    12 = 3
    13 = 12
}

# Snapshot the current values for every affected row/column BEFORE any
# writes happen, so that the permutation doesn't clobber data we still
# need to read for a later destination.
$snapshot = @{}
foreach ($srcRow in $rowMap.Keys) {
    foreach ($col in $cols) {
        $addr = "$col$srcRow"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

# Now write the snapshotted values into their destination rows.
foreach ($srcRow in $rowMap.Keys) {
    $dstRow = $rowMap[$srcRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $dstAddr = "$col$dstRow"
        $ws.Range($dstAddr).Value2 = $snapshot[$srcAddr]
    }
}
